$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.668.32"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "1.635.44"
$ws.Range("E3").Value = "  +1.90%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.494"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.252"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0624"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0839"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.16%  "
$ws.Range("D12").Value = "1.865.43"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").Value = "1.646.17"
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.527"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.68%  "
$ws.Range("D16").Value = "26.679.56"
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.60%  "
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.30"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.35%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.93%  "
$ws.Range("E24").Value = "  +1.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.68%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.120"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0521"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.13%  "
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.57%  "
$ws.Range("E34").Value = "  +0.72%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("D36").Value = "1.169.30"
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0167"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.809"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.07%  "
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.505"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("E41").Value = "  +0.41%  "
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.38"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("D44").Value = "1.777.29"
$ws.Range("E44").Value = "  +2.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.18"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.13%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.06%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -2.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.75%  "
$ws.Range("E50").Value = "  +4.69%  "
$ws.Range("E51").Value = "  +0.59%  "
